$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Confirmations")

# --- Header text updates (row 1) ---
$ws.Range("D1").Value = 'Air Force, nominations'
$ws.Range("E1").Value = 'Air Force, nominations, carryover'
$ws.Range("K1").Value = 'Army, nominations'
$ws.Range("L1").Value = 'Army, nominations, carryover'
$ws.Range("R1").Value = 'Civilian, nominations'
$ws.Range("S1").Value = 'Civilian, nominations, carryover'
$ws.Range("AA1").Value = 'Marine Corps, nominations'
$ws.Range("AB1").Value = 'Marine Corps, nominations, carryover'
$ws.Range("AG1").Value = 'Navy, nominations'
$ws.Range("AH1").Value = 'Navy, nominations, carryover'
$ws.Range("AM1").Value = 'Space Force, nominations'
$ws.Range("AN1").Value = 'Space Force, nominations, carryover'
$ws.Range("AR1").Value = 'Total, failed'
$ws.Range("AS1").Value = 'Total, returned'
$ws.Range("AT1").Value = 'Total, confirmed'
$ws.Range("AU1").Value = 'Total, recess reappointment'
$ws.Range("AV1").Value = 'Total, rejected'
$ws.Range("AW1").Value = 'Total, unconfirmed'
$ws.Range("AX1").Value = 'Total, withdrawn'
$ws.Range("AY1").Value = 'Total, nominations'
$ws.Range("AZ1").Value = 'Total, nominations, carryover'

# --- Numeric value updates ---
$ws.Range("D3").Value = 11818
$ws.Range("D5").Value = 12246
$ws.Range("K5").Value = 9918
$ws.Range("R5").Value = 2046
$ws.Range("AA5").Value = 2132
$ws.Range("D7").Value = 12349
$ws.Range("K7").Value = 10500
$ws.Range("R7").Value = 2670
$ws.Range("AA7").Value = 1586
$ws.Range("AG7").Value = 10159
$ws.Range("D9").Value = 13373
$ws.Range("K9").Value = 15052
$ws.Range("R9").Value = 2364
$ws.Range("AA9").Value = 2775
$ws.Range("AG9").Value = 11370
$ws.Range("D11").Value = 6708
$ws.Range("K11").Value = 11201
$ws.Range("R11").Value = 2661
$ws.Range("AG11").Value = 8001
$ws.Range("D13").Value = 9124
$ws.Range("K13").Value = 10278
$ws.Range("R13").Value = 3029
$ws.Range("AA13").Value = 1359
$ws.Range("AG13").Value = 11255
$ws.Range("D15").Value = 6213
$ws.Range("K15").Value = 8720
$ws.Range("R15").Value = 1781
$ws.Range("AA15").Value = 2332
$ws.Range("AG15").Value = 7165
$ws.Range("D17").Value = 6070
$ws.Range("K17").Value = 5479
$ws.Range("R17").Value = 1782
$ws.Range("AG17").Value = 5047
$ws.Range("D19").Value = 5769
$ws.Range("K19").Value = 6401
$ws.Range("R19").Value = 1581
$ws.Range("AA19").Value = 2826
$ws.Range("AG19").Value = 5585
$ws.Range("D21").Value = 5813
$ws.Range("K21").Value = 6182
$ws.Range("R21").Value = 2028
$ws.Range("AA21").Value = 2978
$ws.Range("D23").Value = 6077
$ws.Range("K23").Value = 5324
$ws.Range("R23").Value = 4417
$ws.Range("AA23").Value = 1227
$ws.Range("AG23").Value = 7375
$ws.Range("D25").Value = 7730
$ws.Range("K25").Value = 9177
$ws.Range("R25").Value = 2956
$ws.Range("AA25").Value = 1291
$ws.Range("AG25").Value = 7015
$ws.Range("D27").Value = 5926
$ws.Range("K27").Value = 6406
$ws.Range("R27").Value = 3444
$ws.Range("AA27").Value = 1565
$ws.Range("AG27").Value = 4749
$ws.Range("D29").Value = 6600
$ws.Range("K29").Value = 7486
$ws.Range("R29").Value = 2675
$ws.Range("AA29").Value = 627
$ws.Range("AG29").Value = 4448
$ws.Range("D31").Value = 6298
$ws.Range("K31").Value = 7300
$ws.Range("R31").Value = 5019
$ws.Range("AG31").Value = 3872
$ws.Range("R33").Value = 4237
$ws.Range("D35").Value = 7387
$ws.Range("K35").Value = 4150
$ws.Range("R35").Value = 2488
$ws.Range("AA35").Value = 1243
$ws.Range("AG35").Value = 4403
$ws.Range("D37").Value = 6198
$ws.Range("K37").Value = 7271
$ws.Range("R37").Value = 2587
$ws.Range("AG37").Value = 4452
$ws.Range("K39").Value = 6386
$ws.Range("R39").Value = 2145
$ws.Range("AG39").Value = 4680
$ws.Range("D41").Value = 6372
$ws.Range("K41").Value = 6242
$ws.Range("R41").Value = 1837
$ws.Range("AA41").Value = 234
$ws.Range("AG41").Value = 4523
$ws.Range("AM41").Value = 470
